$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray "NA" text in C132. It becomes an empty TEXT cell (not a
# fully blank cell), matching every other "no page number" row in the
# sheet (e.g. C2..C132 siblings are empty inline strings, not truly
# blank). A lone leading apostrophe forces Excel to store "" as text
# instead of clearing the cell to a blank/number, then the style is reset
# back to the sheet's normal (unstyled) cell style.
$ws.Range("C132").Value = "'"
$ws.Range("C132").Style = "Normal"

# Append the three new result rows produced by the latest script run.
# The date column must stay plain text ("2025-06-16"), matching every
# other row in the sheet, rather than Excel's automatic date serial
# conversion - so format the cell as Text before typing the value, then
# drop back to the sheet's normal (unstyled) cell style.
function Set-TextValue($cell, $text) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue "A133" "2025-06-16"
$ws.Range("B133").Value = "bonnes pratiques"
$ws.Range("C133").Value = 130
$ws.Range("D133").Value = 1

Set-TextValue "A134" "2025-06-16"
$ws.Range("B134").Value = "bonnes pratiques"
$ws.Range("C134").Value = 135
$ws.Range("D134").Value = 1

Set-TextValue "A135" "2025-06-16"
$ws.Range("B135").Value = "bonnes pratiques"
$ws.Range("C135").Value = 139
$ws.Range("D135").Value = 1
